# Auto-generated update script
# Commit message: Update automàtic: dades i banners [2026-02-22 17:50]
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These cells hold plain-text percentages (e.g. "61%"); force Text
# format first so Excel does not reinterpret them as numeric percents.
foreach ($addr in @("H4", "H13", "H17", "H18", "H24", "H25", "H28", "H30", "H31", "H32", "H35", "H43", "H45", "H46")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("E2").Value = "2026-02-22 17:48:17"
$ws.Range("K2").Value = "13.0 MJ/m2"
$ws.Range("E3").Value = "2026-02-22 17:48:19"
$ws.Range("K3").Value = "15.7 MJ/m2"
$ws.Range("E4").Value = "2026-02-22 17:48:22"
$ws.Range("H4").Value = "61%"
$ws.Range("E5").Value = "2026-02-22 17:48:24"
$ws.Range("E6").Value = "2026-02-22 17:48:26"
$ws.Range("O6").Value = "13.0 °C"
$ws.Range("E7").Value = "2026-02-22 17:48:28"
$ws.Range("J7").Value = "1027.5 hPa"
$ws.Range("O7").Value = "14.3 °C"
$ws.Range("E8").Value = "2026-02-22 17:48:31"
$ws.Range("J8").Value = "1027.0 hPa"
$ws.Range("K8").Value = "16.1 MJ/m2"
$ws.Range("E9").Value = "2026-02-22 17:48:33"
$ws.Range("E10").Value = "2026-02-22 17:48:36"
$ws.Range("E11").Value = "2026-02-22 17:48:38"
$ws.Range("O11").Value = "8.7 °C"
$ws.Range("E12").Value = "2026-02-22 17:48:40"
$ws.Range("O12").Value = "9.8 °C"
$ws.Range("E13").Value = "2026-02-22 17:48:42"
$ws.Range("H13").Value = "62%"
$ws.Range("J13").Value = "1030.7 hPa"
$ws.Range("K13").Value = "15.7 MJ/m2"
$ws.Range("O13").Value = "6.1 °C"
$ws.Range("E14").Value = "2026-02-22 17:48:45"
$ws.Range("E15").Value = "2026-02-22 17:48:47"
$ws.Range("O15").Value = "11.1 °C"
$ws.Range("E16").Value = "2026-02-22 17:48:49"
$ws.Range("K16").Value = "13.1 MJ/m2"
$ws.Range("O16").Value = "5.5 °C"
$ws.Range("E17").Value = "2026-02-22 17:48:51"
$ws.Range("H17").Value = "27%"
$ws.Range("E18").Value = "2026-02-22 17:48:54"
$ws.Range("H18").Value = "74%"
$ws.Range("J18").Value = "1027.7 hPa"
$ws.Range("K18").Value = "15.3 MJ/m2"
$ws.Range("O18").Value = "9.9 °C"
$ws.Range("E19").Value = "2026-02-22 17:48:56"
$ws.Range("E20").Value = "2026-02-22 17:48:59"
$ws.Range("K20").Value = "16.3 MJ/m2"
$ws.Range("E21").Value = "2026-02-22 17:49:01"
$ws.Range("J21").Value = "1029.4 hPa"
$ws.Range("K21").Value = "15.5 MJ/m2"
$ws.Range("O21").Value = "8.9 °C"
$ws.Range("E22").Value = "2026-02-22 17:49:03"
$ws.Range("K22").Value = "16.3 MJ/m2"
$ws.Range("E23").Value = "2026-02-22 17:49:06"
$ws.Range("K23").Value = "15.8 MJ/m2"
$ws.Range("E24").Value = "2026-02-22 17:49:08"
$ws.Range("H24").Value = "82%"
$ws.Range("J24").Value = "1029.9 hPa"
$ws.Range("K24").Value = "15.6 MJ/m2"
$ws.Range("O24").Value = "7.5 °C"
$ws.Range("E25").Value = "2026-02-22 17:49:11"
$ws.Range("H25").Value = "26%"
$ws.Range("O25").Value = "7.4 °C"
$ws.Range("E26").Value = "2026-02-22 17:49:13"
$ws.Range("O26").Value = "11.8 °C"
$ws.Range("E27").Value = "2026-02-22 17:49:15"
$ws.Range("K27").Value = "16.3 MJ/m2"
$ws.Range("E28").Value = "2026-02-22 17:49:18"
$ws.Range("H28").Value = "63%"
$ws.Range("J28").Value = "1027.8 hPa"
$ws.Range("K28").Value = "15.1 MJ/m2"
$ws.Range("O28").Value = "10.5 °C"
$ws.Range("E29").Value = "2026-02-22 17:49:20"
$ws.Range("O29").Value = "10.1 °C"
$ws.Range("E30").Value = "2026-02-22 17:49:23"
$ws.Range("H30").Value = "69%"
$ws.Range("O30").Value = "12.6 °C"
$ws.Range("E31").Value = "2026-02-22 17:49:25"
$ws.Range("H31").Value = "64%"
$ws.Range("J31").Value = "1026.7 hPa"
$ws.Range("O31").Value = "14.2 °C"
$ws.Range("E32").Value = "2026-02-22 17:49:27"
$ws.Range("H32").Value = "68%"
$ws.Range("K32").Value = "16.0 MJ/m2"
$ws.Range("O32").Value = "6.4 °C"
$ws.Range("E33").Value = "2026-02-22 17:49:30"
$ws.Range("J33").Value = "1028.9 hPa"
$ws.Range("K33").Value = "15.7 MJ/m2"
$ws.Range("L33").Value = "13.0 km/h - 98º 17:02 TU"
$ws.Range("O33").Value = "8.1 °C"
$ws.Range("E34").Value = "2026-02-22 17:49:32"
$ws.Range("O34").Value = "4.7 °C"
$ws.Range("E35").Value = "2026-02-22 17:49:35"
$ws.Range("H35").Value = "40%"
$ws.Range("K35").Value = "16.3 MJ/m2"
$ws.Range("E36").Value = "2026-02-22 17:49:37"
$ws.Range("J36").Value = "1027.5 hPa"
$ws.Range("K36").Value = "15.2 MJ/m2"
$ws.Range("O36").Value = "11.7 °C"
$ws.Range("E37").Value = "2026-02-22 17:49:39"
$ws.Range("J37").Value = "1029.9 hPa"
$ws.Range("O37").Value = "8.1 °C"
$ws.Range("E38").Value = "2026-02-22 17:49:42"
$ws.Range("O38").Value = "11.5 °C"
$ws.Range("E39").Value = "2026-02-22 17:49:44"
$ws.Range("E40").Value = "2026-02-22 17:49:47"
$ws.Range("J40").Value = "1029.3 hPa"
$ws.Range("O40").Value = "10.3 °C"
$ws.Range("E41").Value = "2026-02-22 17:49:49"
$ws.Range("K41").Value = "15.6 MJ/m2"
$ws.Range("E42").Value = "2026-02-22 17:49:51"
$ws.Range("O42").Value = "10.7 °C"
$ws.Range("E43").Value = "2026-02-22 17:49:54"
$ws.Range("H43").Value = "72%"
$ws.Range("K43").Value = "15.3 MJ/m2"
$ws.Range("O43").Value = "8.9 °C"
$ws.Range("E44").Value = "2026-02-22 17:49:56"
$ws.Range("K44").Value = "15.7 MJ/m2"
$ws.Range("E45").Value = "2026-02-22 17:49:59"
$ws.Range("H45").Value = "52%"
$ws.Range("J45").Value = "1028.8 hPa"
$ws.Range("K45").Value = "13.9 MJ/m2"
$ws.Range("E46").Value = "2026-02-22 17:50:01"
$ws.Range("H46").Value = "75%"
$ws.Range("J46").Value = "1029.9 hPa"
$ws.Range("O46").Value = "8.8 °C"
